$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 16: new journal entry (16/05/2018)
$ws.Range("A16").Value = (Get-Date -Year 2018 -Month 5 -Day 16 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B16").Value = 0.5
$ws.Range("C16").Value = "Finalisation des use-case & scénarios"
$ws.Range("E16").Value = 9

# Row 17: new journal entry (continuation, no date)
$ws.Range("B17").Value = 3
$ws.Range("C17").Value = "Création de l'application de base"
$ws.Range("D17").Value = "Mise en place d'un fichier de config, création d'une classe ""Main_Window"", et création d'un semblant d'interface grâce à tkinter"
$ws.Rows.Item(17).RowHeight = 30

# Update the active selection to match the new edit location
$ws.Range("B18").Select()
